$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two more mismatched example rows beneath the existing ones.
$ws.Range("A3").Value = "Instructor"
$ws.Range("B3").Value = "John"
$ws.Range("C3").Value = "Distribution Files\S1984\MecE_265_Aryanci\MecE_265.xlsx"

$ws.Range("A4").Value = "Course Name"
$ws.Range("B4").Value = "MATH 101"
$ws.Range("C4").Value = "Distribution Files\W2009\MEC E 260 - 502 (77588) Ayranci\MecE260_UNDERGRAD Grade Dist Form.xlsx"

# Widen column C so the longer file-path strings are readable (target ~82.4765625 chars;
# the host snaps ColumnWidth to its internal pixel grid, so feed the value that lands
# closest to the target after rounding).
$ws.Columns.Item(3).ColumnWidth = 81.667
